$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 125
$ws.Cells.Item(8, 8).Value = 92.5
$ws.Cells.Item(8, 9).Value = 92.5
$ws.Cells.Item(8, 11).Value = 277.5
$ws.Cells.Item(8, 13).Value = -138.5
$ws.Cells.Item(17, 8).Value = 2882.2942
$ws.Cells.Item(17, 10).Value = 2599.8
$ws.Cells.Item(17, 12).Value = 7799.400000000001
$ws.Cells.Item(17, 14).Value = -8135.400000000001
$ws.Cells.Item(18, 8).Value = 592.2727
$ws.Cells.Item(18, 9).Value = 601.5
$ws.Cells.Item(18, 11).Value = 601.5
$ws.Cells.Item(18, 13).Value = -317.5
$ws.Cells.Item(19, 8).Value = 998.5
$ws.Cells.Item(19, 9).Value = 998.5
$ws.Cells.Item(19, 11).Value = 998.5
$ws.Cells.Item(19, 13).Value = -823.5
$ws.Cells.Item(32, 8).Value = 4333.3335
$ws.Cells.Item(32, 9).Value = 4700
$ws.Cells.Item(32, 10).Value = 4150
$ws.Cells.Item(32, 11).Value = 4700
$ws.Cells.Item(32, 12).Value = 4150
$ws.Cells.Item(32, 13).Value = -4374
$ws.Cells.Item(32, 14).Value = -4802
$ws.Cells.Item(39, 8).Value = 204.6923
$ws.Cells.Item(39, 9).Value = 116.2
$ws.Cells.Item(39, 11).Value = 348.6
$ws.Cells.Item(39, 13).Value = -52.60000000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1008.7143
$ws.Cells.Item(2, 9).Value = 1150.6666
$ws.Cells.Item(2, 10).Value = 157
$ws.Cells.Item(2, 11).Value = 1150.6666
$ws.Cells.Item(2, 12).Value = 157
$ws.Cells.Item(2, 13).Value = -1037.6666
$ws.Cells.Item(2, 14).Value = -383
$ws.Cells.Item(116, 8).Value = 1008.7143
$ws.Cells.Item(116, 9).Value = 1150.6666
$ws.Cells.Item(116, 11).Value = 1150.6666
$ws.Cells.Item(116, 13).Value = 1143.3334
$ws.Cells.Item(116, 14).Value = -4745
$ws.Cells.Item(132, 8).Value = 3470.6667
$ws.Cells.Item(132, 9).Value = 3470.6667
$ws.Cells.Item(132, 11).Value = 10412.0001
$ws.Cells.Item(132, 13).Value = -7882.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1008.7143
$ws.Cells.Item(3, 9).Value = 1150.6666
$ws.Cells.Item(3, 10).Value = 157
$ws.Cells.Item(3, 11).Value = 1150.6666
$ws.Cells.Item(3, 12).Value = 157
$ws.Cells.Item(3, 13).Value = -1036.6666
$ws.Cells.Item(3, 14).Value = -385
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 13).ClearContents() | Out-Null
$ws.Cells.Item(40, 8).Value = 59999
$ws.Cells.Item(40, 10).Value = 59999
$ws.Cells.Item(40, 12).Value = 59999
$ws.Cells.Item(40, 14).Value = -60529
$ws.Cells.Item(80, 8).Value = 893.4
$ws.Cells.Item(80, 9).Value = 890.5
$ws.Cells.Item(80, 10).Value = 895.3333
$ws.Cells.Item(80, 11).Value = 890.5
$ws.Cells.Item(80, 12).Value = 895.3333
$ws.Cells.Item(80, 13).Value = 107.5
$ws.Cells.Item(80, 14).Value = -2891.3333
$ws.Cells.Item(83, 8).Value = 893.4
$ws.Cells.Item(83, 9).Value = 890.5
$ws.Cells.Item(83, 10).Value = 895.3333
$ws.Cells.Item(83, 11).Value = 4452.5
$ws.Cells.Item(83, 12).Value = 4476.6665
$ws.Cells.Item(83, 13).Value = 539.5
$ws.Cells.Item(83, 14).Value = -14460.6665
$ws.Cells.Item(94, 8).Value = 2500
$ws.Cells.Item(94, 9).Value = 2000
$ws.Cells.Item(94, 10).Value = 3000
$ws.Cells.Item(94, 11).Value = 2000
$ws.Cells.Item(94, 12).Value = 3000
$ws.Cells.Item(94, 13).Value = -1549
$ws.Cells.Item(94, 14).Value = -3902
$ws.Cells.Item(134, 8).Value = 4125
$ws.Cells.Item(134, 9).Value = 4125
$ws.Cells.Item(134, 11).Value = 12375
$ws.Cells.Item(134, 13).Value = -9840
$ws.Cells.Item(135, 8).Value = 50000
$ws.Cells.Item(135, 10).Value = 50000
$ws.Cells.Item(135, 12).Value = 50000
$ws.Cells.Item(135, 14).Value = -60140

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 634
$ws.Cells.Item(109, 9).Value = 765.3333
$ws.Cells.Item(109, 10).Value = 240
$ws.Cells.Item(109, 11).Value = 2295.9999
$ws.Cells.Item(109, 12).Value = 720
$ws.Cells.Item(109, 13).Value = -1255.9999
$ws.Cells.Item(109, 14).Value = -2800
$ws.Cells.Item(129, 8).Value = 1331
$ws.Cells.Item(129, 9).Value = 563
$ws.Cells.Item(129, 11).Value = 1689
$ws.Cells.Item(129, 13).Value = 3311
$ws.Cells.Item(141, 8).Value = 1900
$ws.Cells.Item(141, 9).Value = 1900
$ws.Cells.Item(141, 11).Value = 5700
$ws.Cells.Item(141, 13).Value = -520

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 7999.5
$ws.Cells.Item(18, 10).Value = 7999.5
$ws.Cells.Item(18, 12).Value = 7999.5
$ws.Cells.Item(18, 14).Value = -8585.5
$ws.Cells.Item(21, 8).Value = 500
$ws.Cells.Item(21, 10).Value = 500
$ws.Cells.Item(21, 12).Value = 500
$ws.Cells.Item(21, 14).Value = -846
$ws.Cells.Item(30, 8).Value = 500
$ws.Cells.Item(30, 10).Value = 500
$ws.Cells.Item(30, 12).Value = 500
$ws.Cells.Item(30, 14).Value = -710
$ws.Cells.Item(44, 8).Value = 25007.75
$ws.Cells.Item(44, 10).Value = 25007.75
$ws.Cells.Item(44, 12).Value = 25007.75
$ws.Cells.Item(44, 14).Value = -26199.75
$ws.Cells.Item(47, 8).Value = 30000
$ws.Cells.Item(47, 10).Value = 30000
$ws.Cells.Item(47, 12).Value = 30000
$ws.Cells.Item(47, 14).Value = -31136
$ws.Cells.Item(70, 8).Value = 2999
$ws.Cells.Item(70, 10).Value = 2999
$ws.Cells.Item(70, 12).Value = 2999
$ws.Cells.Item(70, 14).Value = -3539
$ws.Cells.Item(73, 8).Value = 2999
$ws.Cells.Item(73, 10).Value = 2999
$ws.Cells.Item(73, 12).Value = 2999
$ws.Cells.Item(73, 14).Value = -4871
$ws.Cells.Item(93, 8).Value = 41000
$ws.Cells.Item(93, 10).Value = 41000
$ws.Cells.Item(93, 12).Value = 41000
$ws.Cells.Item(93, 14).Value = -44744
$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 12).Value = 0
$ws.Cells.Item(117, 14).ClearContents() | Out-Null
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).ClearContents() | Out-Null
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(130, 14).ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(5, 8).Value = 5000
$ws.Cells.Item(5, 9).Value = 5000
$ws.Cells.Item(5, 11).Value = 5000
$ws.Cells.Item(5, 13).Value = -4887
$ws.Cells.Item(22, 8).Value = 999.75
$ws.Cells.Item(22, 9).Value = 999.5
$ws.Cells.Item(22, 11).Value = 999.5
$ws.Cells.Item(22, 13).Value = -704.5
$ws.Cells.Item(27, 8).Value = 999.75
$ws.Cells.Item(27, 9).Value = 999.5
$ws.Cells.Item(27, 11).Value = 999.5
$ws.Cells.Item(27, 13).Value = -892.5
$ws.Cells.Item(43, 8).Value = 5000
$ws.Cells.Item(43, 10).Value = 5000
$ws.Cells.Item(43, 12).Value = 5000
$ws.Cells.Item(43, 14).Value = -5386
$ws.Cells.Item(61, 8).Value = 1442.5
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 13).ClearContents() | Out-Null
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 13).ClearContents() | Out-Null
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 13).ClearContents() | Out-Null
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents() | Out-Null
$ws.Cells.Item(113, 8).Value = 1442.5
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 13).ClearContents() | Out-Null
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).ClearContents() | Out-Null
$ws.Cells.Item(136, 8).Value = 2389888.2
$ws.Cells.Item(136, 9).Value = 1667999.4
$ws.Cells.Item(136, 11).Value = 5003998.199999999
$ws.Cells.Item(136, 13).Value = -5001448.199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 9).Value = 3000
$ws.Cells.Item(18, 10).Value = 2997
$ws.Cells.Item(18, 11).Value = 3000
$ws.Cells.Item(18, 12).Value = 2997
$ws.Cells.Item(18, 13).Value = -2827
$ws.Cells.Item(18, 14).Value = -3343
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 14).ClearContents() | Out-Null
$ws.Cells.Item(124, 8).Value = 49999.5
$ws.Cells.Item(124, 9).Value = 49999
$ws.Cells.Item(124, 11).Value = 49999
$ws.Cells.Item(124, 13).Value = -45089
$ws.Cells.Item(136, 8).Value = 1915.2174
$ws.Cells.Item(136, 9).Value = 1852.5
$ws.Cells.Item(136, 10).Value = 2333.3333
$ws.Cells.Item(136, 11).Value = 5557.5
$ws.Cells.Item(136, 12).Value = 6999.999899999999
$ws.Cells.Item(136, 13).Value = -3007.5
$ws.Cells.Item(136, 14).Value = -12099.9999
